# Fruta / hortaliza, semanal
# Insert a new weekly price-report row for "Zapallo / Camote" right after the
# existing row 27, pushing all the following rows down by one (rows 28-52
# become 29-53). Column D (Fecha) keeps its "s=2" date style on the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 28 - shifts existing rows 28:52 down to 29:53
$ws.Rows(28).Insert()

# Populate the newly inserted row 28 with the latest week's data
$ws.Range("A28").Value = 1
$ws.Range("B28").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C28").Value = "Arica y Parinacota"
$ws.Range("D28").Value = 44942
$ws.Range("E28").Value = 15
$ws.Range("F28").Value = 100112045
$ws.Range("G28").Value = "Zapallo"
$ws.Range("H28").Value = "Camote"
$ws.Range("I28").Value = "1a nueva(o)"
$ws.Range("J28").Value = 450
$ws.Range("K28").Value = 430
$ws.Range("L28").Value = 450
$ws.Range("M28").Value = 441
$ws.Range("N28").Value = "$/kilo (volumen en unidades)"
$ws.Range("O28").Value = "Perú"
$ws.Range("P28").Value = 441
$ws.Range("Q28").Value = 1
$ws.Range("R28").Value = "Hortaliza"
